# Update "Sorter Inspection Validation Atlanta.xlsx"
# 1) Inspection Log sheet: widen column A and append a new inspection row (204)
# 2) Weekly Summary sheet: append a repeated header row (23) and a new data row (24)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Inspection Log"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Inspection Log")

# Widen column A from 12 to 21 (character width as stored in the OOXML).
# ColumnWidth undergoes a pixel round-trip conversion in Excel, so 20.17 is the
# value that yields a stored width of exactly 21.
$ws1.Range("A1").EntireColumn.ColumnWidth = 20.17

# Append new row 204 with the same look/formatting as row 203 (date format on A,
# green Pass fill on C and E), then set its values.
$ws1.Range("A203:E203").Copy()
$ws1.Range("A204:E204").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Cells.Item(204, 1).Value = 45798
$ws1.Cells.Item(204, 2).Value = "Strand 7"
$ws1.Cells.Item(204, 3).Value = "1:16:09"
$ws1.Cells.Item(204, 4).Value = "No"
$ws1.Cells.Item(204, 5).Value = "Pass"

# ---------------------------------------------------------------------------
# Sheet 2: "Weekly Summary"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Weekly Summary")

# Row 23: repeat of the header labels (no special formatting).
$ws2.Cells.Item(23, 1).Value = "Week Range"
$ws2.Cells.Item(23, 2).Value = "Strands Completed"
$ws2.Cells.Item(23, 3).Value = "All 8 Present"

# Row 24: new weekly summary entry, formatted like row 22 (red Fail fill on C).
$ws2.Range("A22:C22").Copy()
$ws2.Range("A24:C24").PasteSpecial(-4122)  # xlPasteFormats

$ws2.Cells.Item(24, 1).Value = "05-18-25 to 05-24-25"
$ws2.Cells.Item(24, 2).Value = "Strand 7"
$ws2.Cells.Item(24, 3).Value = "Fail"
